$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently has the teacher-name / dates mis-aligned against their
# row labels (rows 10-21): e.g. "Objetivos:" (A10) sits next to the teacher
# name instead of the PT objectives text, "Programa resumido:" (A13) sits next
# to "Semestral" instead of a short syllabus, etc. Insert a blank row at 13 so
# everything from the old row 13 down shifts to row 14 down (matching how the
# corrected sheet re-aligns labels with their real content), then rewrite the
# cell contents/heights for rows 10-22 to the corrected & expanded text.
$ws.Rows.Item(13).Insert()

# Row 13 is now a brand-new blank row; give its B/C cells the same formatting
# (wrap text, vertical-top, plain vs. red font) as the other data columns use,
# then drop the placeholder A13 cell the Insert left behind (row 13 has no A
# entry in the corrected sheet).
$ws.Range("B10").Copy()
$ws.Range("B13").PasteSpecial(-4122)
$ws.Range("C10").Copy()
$ws.Range("C13").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("A13").Clear()

# --- Row 10 ---
$ws.Range('A10').Value = 'Objetivos:'
$ws.Range('B10').Value = 'Familiarizar o aluno com a utilização de instrumentos de medidas mecânicas. Elaboração de tabelas e gráficos com escalas lineares e logarítmicas. Introdução de conceitos básicos da teoria de Erros e do Método dos Mínimos Quadrados. Realização de experimentos básicos de mecânica e elaboração de relatórios.'
$ws.Range('C10').Value = 'Familiarizar o aluno com a utilização de instrumentos de medidas mecânicas. Elaboração de tabelas e gráficos com escalas lineares e logarítmicas. Introdução de conceitos básicos da teoria de Erros e do Método dos Mínimos Quadrados. Realização de experimentos básicos de mecânica e elaboração de relatórios.'
$ws.Rows.Item(10).RowHeight = 60

# --- Row 11 ---
$ws.Range('A11').Value = 'Objectives:'
$ws.Range('B11').Value = 'To familiarize the student with the use of measuring instruments. Drafting tables and graphics with linear and logarithmic scales . Basic Concepts of Error Theory and method of least squares. Basic mechanics experiments and preparation of reports.'
$ws.Range('C11').Value = 'To familiarize the student with the use of measuring instruments. Drafting tables and graphics with linear and logarithmic scales . Basic Concepts of Error Theory and method of least squares. Basic mechanics experiments and preparation of reports.'
$ws.Rows.Item(11).RowHeight = 60

# --- Row 12 ---
$ws.Range('A12').Value = 'Docentes responsáveis:'
$ws.Range('B12').ClearContents()
$ws.Range('C12').ClearContents()
$ws.Rows.Item(12).AutoFit()

# --- Row 13 ---
$ws.Range('A13').ClearContents()
$ws.Range('B13').Value = '9149242 - Fernando Catalani'
$ws.Range('C13').Value = '9149242 - Fernando Catalani'
$ws.Rows.Item(13).AutoFit()

# --- Row 14 ---
$ws.Range('A14').Value = 'Programa resumido:'
$ws.Range('B14').Value = 'Instrumentos de medidas; Construção de Tabelas e Gráficos; Método dos mínimos quadrados; Estática, Cinemática; Dinâmica; Conservação de Energia Mecânica; Choques Unidimensionais'
$ws.Range('C14').Value = 'Instrumentos de medidas; Construção de Tabelas e Gráficos; Método dos mínimos quadrados; Estática, Cinemática; Dinâmica; Conservação de Energia Mecânica; Choques Unidimensionais'
$ws.Rows.Item(14).RowHeight = 60

# --- Row 15 ---
$ws.Range('A15').Value = 'Short syllabus:'
$ws.Range('B15').ClearContents()
$ws.Range('C15').ClearContents()
$ws.Rows.Item(15).RowHeight = 60

# --- Row 16 ---
$ws.Range('A16').Value = 'Programa:'
$ws.Range('B16').Value = '1) Instrumentos de medidas. Estimativa de erro nas medidas, propagação de erros e algarismos significativos.2) Construção de Tabelas e Gráficos. Linearização.3) Regressão linear. Introdução ao método dos mínimos quadrados. 4) Cinemática. Movimento Retilíneo Uniforme e Movimento Retilíneo uniformemente variado. Queda Livre.5) Estática. Equilíbrio de um ponto Material. 6) Atrito.7) Lei de Hooke. Módulo de Young. 8) Conservação de Energia. Conceito de Conservação da Energia Mecânica. Sistema Massa-mola.9) Choques Unidimensionais.'
$ws.Range('C16').Value = '1) Instrumentos de medidas. Estimativa de erro nas medidas, propagação de erros e algarismos significativos.2) Construção de Tabelas e Gráficos. Linearização.3) Regressão linear. Introdução ao método dos mínimos quadrados. 4) Cinemática. Movimento Retilíneo Uniforme e Movimento Retilíneo uniformemente variado. Queda Livre.5) Estática. Equilíbrio de um ponto Material. 6) Atrito.7) Lei de Hooke. Módulo de Young. 8) Conservação de Energia. Conceito de Conservação da Energia Mecânica. Sistema Massa-mola.9) Choques Unidimensionais.'
$ws.Rows.Item(16).RowHeight = 120

# --- Row 17 ---
$ws.Range('A17').Value = 'Syllabus:'
$ws.Range('B17').Value = '1) Simple measures. Error Estimation of measures. Error propagation and significant figures.2) Construction of Tables and Graphs. Linearization.3) Introduction to the method of squares linear regression minimum.4) Kinematics. Rectilinear motion and uniformly varied motion. Free fall.5) Statics. Equilibrium of a material point. 6) Friction.7) Hooke''s Law. Young´s Modulus.8) Energy conservation. Conservation Concept of Energy Mechanics. Mass-spring system.9) Shocks.'
$ws.Range('C17').Value = '1) Simple measures. Error Estimation of measures. Error propagation and significant figures.2) Construction of Tables and Graphs. Linearization.3) Introduction to the method of squares linear regression minimum.4) Kinematics. Rectilinear motion and uniformly varied motion. Free fall.5) Statics. Equilibrium of a material point. 6) Friction.7) Hooke''s Law. Young´s Modulus.8) Energy conservation. Conservation Concept of Energy Mechanics. Mass-spring system.9) Shocks.'
$ws.Rows.Item(17).RowHeight = 120

# --- Row 18 ---
$ws.Range('A18').Value = 'Avaliação:'
$ws.Range('B18').ClearContents()
$ws.Range('C18').ClearContents()
$ws.Rows.Item(18).AutoFit()

# --- Row 19 ---
$ws.Range('A19').Value = 'Método:'
$ws.Range('B19').Value = 'NF=A avaliação será composta por provas, listas, projetos, seminários e outras formas que farão a composição das notas, sendo estipulada a média final a somatória destas notas (N), com no mínimo duas avaliações, sendo: (N1+...+Nn)/n.'
$ws.Range('C19').Value = 'NF=A avaliação será composta por provas, listas, projetos, seminários e outras formas que farão a composição das notas, sendo estipulada a média final a somatória destas notas (N), com no mínimo duas avaliações, sendo: (N1+...+Nn)/n.'
$ws.Rows.Item(19).RowHeight = 60

# --- Row 20 ---
$ws.Range('A20').Value = 'Critério:'
$ws.Range('B20').Value = 'NF≥ 5,0.'
$ws.Range('C20').Value = 'NF≥ 5,0.'
$ws.Rows.Item(20).RowHeight = 60

# --- Row 21 ---
$ws.Range('A21').Value = 'Norma de recuperação:'
$ws.Range('B21').Value = 'O (NF+RC)/2 ≥ 5,0, onde RC é uma prova de recuperação a ser aplicada.'
$ws.Range('C21').Value = 'O (NF+RC)/2 ≥ 5,0, onde RC é uma prova de recuperação a ser aplicada.'
$ws.Rows.Item(21).RowHeight = 60

# --- Row 22 ---
$ws.Range('A22').Value = 'Bibliografia:'
$ws.Range('B22').Value = 'Apostilas do Laboratório de Ensino de Física do IFSC/USP.CRUZ, C. H. B.; FRAGNITO, H. L.; COSTA, I. F.; MELLO, B. A. Guia do Curso deLaboratório: Física Experimental I, IFGW/UNICAMP (2005).NUSSENZVEIG, H.M. Curso de Física Básica. Vol. 1, Edgard Blucher (2008).RESNICK, R.; HALLIDAY, D. Fundamentos de Física. Vol.1, LTC (2008).TIPLER, P.; MOSCA, G. Física para Cientistas e Engenheiros. Vol.1, LTC (2008).SEARS, F. W.; ZEMANSKY, M. W.; YOUNG, H. D.; FREEDMAN, R. A. Física I, Vol. 1, Pearson Addison Wesley (2009).JEWETT Jr, John W.; SERWAY, Raymond A. Princípios de Física. Vol. 1, Thomson Pioneira (2008).'
$ws.Range('C22').Value = 'Apostilas do Laboratório de Ensino de Física do IFSC/USP.CRUZ, C. H. B.; FRAGNITO, H. L.; COSTA, I. F.; MELLO, B. A. Guia do Curso deLaboratório: Física Experimental I, IFGW/UNICAMP (2005).NUSSENZVEIG, H.M. Curso de Física Básica. Vol. 1, Edgard Blucher (2008).RESNICK, R.; HALLIDAY, D. Fundamentos de Física. Vol.1, LTC (2008).TIPLER, P.; MOSCA, G. Física para Cientistas e Engenheiros. Vol.1, LTC (2008).SEARS, F. W.; ZEMANSKY, M. W.; YOUNG, H. D.; FREEDMAN, R. A. Física I, Vol. 1, Pearson Addison Wesley (2009).JEWETT Jr, John W.; SERWAY, Raymond A. Princípios de Física. Vol. 1, Thomson Pioneira (2008).'
$ws.Rows.Item(22).RowHeight = 120

